$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Participants" header in F1, matching the style of the other headers
$ws.Range("F1").Value = "Participants"
$ws.Range("E1:E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Fill in participant counts per cohort (100s=5, 200s=1, 400s=15), 4 rows each
$participants = @(5, 5, 5, 5, 1, 1, 1, 1, 15, 15, 15, 15)
for ($i = 0; $i -lt $participants.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $participants[$i]
}
